# edit.ps1 - applies categorized_data.xlsx changes per commit diff
# Adds new transaction rows to Income/Expenses/Business Expenses/Subscriptions/
# Uncertain Expenses sheets, updates Weekly Budget + Balance Summary aggregates,
# and fixes the Balances sheet bank total.

$wb = $excel.ActiveWorkbook

function Set-TxnRow($ws, $r, $date, $amt, $desc, $src) {
    # Column A holds text dates like "12/20/2023" - force text so Excel
    # does not auto-convert them into date serial numbers.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $date
    $ws.Cells.Item($r, 1).Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $amt
    $ws.Cells.Item($r, 3).Value = $desc
    $ws.Cells.Item($r, 4).Value = $src
}

# --- Income: append new transaction rows ---
$ws = $wb.Worksheets.Item("Income")
Set-TxnRow $ws 180 "12/20/2023" 104 "PATRICK RIVERA CO: PATRICK RIVERA : IAT PAYPAL  ID: 770510487C DATA: FF3               US %% ACH ECC IAT" "Employer"
Set-TxnRow $ws 181 "12/15/2023" 3000 "TIJO VARUGHESE : P2P PAYMNT  ID: 9000041902CO: American Heritag  NAME: TIJO VARUGHESE %% ACH ECC WEB" "Upwork"
Set-TxnRow $ws 182 "11/30/2023" 1000 "TIJO VARUGHESE : P2P PAYMNT  ID: 9000041902CO: American Heritag  NAME: TIJO VARUGHESE %% ACH ECC WEB" "Upwork"
Set-TxnRow $ws 183 "11/27/2023" 300 "TIJO VARUGHESE : P2P PAYMNT  ID: 9000041902CO: American Heritag  NAME: TIJO VARUGHESE %% ACH ECC WEB" "Upwork"
Set-TxnRow $ws 184 "11/22/2023" 104 "PATRICK RIVERA CO: PATRICK RIVERA : IAT PAYPAL  ID: 770510487C DATA: FF3               US %% ACH ECC IAT" "Employer"
Set-TxnRow $ws 185 "11/22/2023" 200 "TIJO VARUGHESE : P2P PAYMNT  ID: 9000041902CO: American Heritag  NAME: TIJO VARUGHESE %% ACH ECC WEB" "Upwork"
Set-TxnRow $ws 186 "10/31/2023" 300 "TIJO VARUGHESE : P2P PAYMNT  ID: 9000041902CO: American Heritag  NAME: TIJO VARUGHESE %% ACH ECC WEB" "American Heritage"
Set-TxnRow $ws 187 "10/20/2023" 100 "From Share 00" "Bank"
Set-TxnRow $ws 188 "10/10/2023" 54 "From Share 00" "Bank"
Set-TxnRow $ws 189 "12/15/2023" 3000 "TIJO VARUGHESE : P2P PAYMNT  ID: 9000041902CO: American Heritag  NAME: TIJO VARUGHESE %% ACH ECC WEB" "Upwork"
Set-TxnRow $ws 190 "11/30/2023" 1000 "TIJO VARUGHESE : P2P PAYMNT  ID: 9000041902CO: American Heritag  NAME: TIJO VARUGHESE %% ACH ECC WEB" "Upwork"
Set-TxnRow $ws 191 "11/29/2023" 210 "From Share 00" "Bank"
Set-TxnRow $ws 192 "11/27/2023" 300 "TIJO VARUGHESE : P2P PAYMNT  ID: 9000041902CO: American Heritag  NAME: TIJO VARUGHESE %% ACH ECC WEB" "Upwork"
Set-TxnRow $ws 193 "11/22/2023" 200 "TIJO VARUGHESE : P2P PAYMNT  ID: 9000041902CO: American Heritag  NAME: TIJO VARUGHESE %% ACH ECC WEB" "Upwork"

# --- Expenses: append new transaction rows ---
$ws = $wb.Worksheets.Item("Expenses")
Set-TxnRow $ws 643 "11/29/2023" 15 "ODP Fee %% ACH ECC WEB %% ACH Trace 091000016840533" "Bank"

# --- Business Expenses: append new transaction rows ---
$ws = $wb.Worksheets.Item("Business Expenses")
Set-TxnRow $ws 276 "12/29/2023" 85 "PAYPAL : INST XFER  ID: PAYPALSI77 DATA: INSTANT TRANSFER  CO: PAYPAL NAME: CREED SERVICES %% ACH ECC WEB" "Business Expense"
Set-TxnRow $ws 277 "12/18/2023" 2861 "PAYPAL : INST XFER  ID: PAYPALSI77 DATA: INSTANT TRANSFER  CO: PAYPAL NAME: CREED SERVICES %% ACH ECC WEB" "Business Expense"
Set-TxnRow $ws 278 "11/29/2023" 954 "PAYPAL : INST XFER  ID: PAYPALSI77 DATA: INSTANT TRANSFER  CO: PAYPAL NAME: CREED SERVICES" "Business Expense"
Set-TxnRow $ws 279 "11/13/2023" 120 "PAYPAL : INST XFER  ID: PAYPALSI77 DATA: INSTANT TRANSFER  CO: PAYPAL NAME: CREED SERVICES %% ACH ECC WEB" "PayPal"
Set-TxnRow $ws 280 "11/02/2023" 105 "PATRICK RIVERA CO: PATRICK RIVERA : IAT PAYPAL  ID: 770510487C DATA: FF3               US %% ACH ECC IAT" "PayPal"
Set-TxnRow $ws 281 "10/20/2023" 100 "PATRICK RIVERA CO: PATRICK RIVERA : IAT PAYPAL  ID: 770510487C DATA: FF3               US %% ACH ECC IAT" "PayPal"
Set-TxnRow $ws 282 "10/10/2023" 105 "PATRICK RIVERA CO: PATRICK RIVERA : IAT PAYPAL  ID: 770510487C DATA: FF3               US %% ACH ECC IAT" "PayPal"
Set-TxnRow $ws 283 "10/10/2023" 35 "PAYPAL : INST XFER  ID: PAYPALSI77 DATA: INSTANT TRANSFER  CO: PAYPAL NAME: CREED SERVICES %% ACH ECC WEB" "PayPal"
Set-TxnRow $ws 284 "09/14/2023" 14 "FIVERR INTERNATI CO: FIVERR INTERNATIONAL LTD. : IAT PAYPAL  ID: 770510487C DATA: FF3               US" "Fiverr"
Set-TxnRow $ws 285 "12/29/2023" 85 "PAYPAL : INST XFER  ID: PAYPALSI77 DATA: INSTANT TRANSFER  CO: PAYPAL NAME: CREED SERVICES %% ACH ECC WEB" "Business Expense"
Set-TxnRow $ws 286 "12/20/2023" 104 "PATRICK RIVERA CO: PATRICK RIVERA : IAT PAYPAL  ID: 770510487C DATA: FF3               US %% ACH ECC IAT" "Employer"
Set-TxnRow $ws 287 "12/18/2023" 2861 "PAYPAL : INST XFER  ID: PAYPALSI77 DATA: INSTANT TRANSFER  CO: PAYPAL NAME: CREED SERVICES %% ACH ECC WEB" "Business Expense"
Set-TxnRow $ws 288 "11/29/2023" 954 "PAYPAL : INST XFER  ID: PAYPALSI77 DATA: INSTANT TRANSFER  CO: PAYPAL NAME: CREED SERVICES" "Business Expense"
Set-TxnRow $ws 289 "11/22/2023" 104 "PATRICK RIVERA CO: PATRICK RIVERA : IAT PAYPAL  ID: 770510487C DATA: FF3               US %% ACH ECC IAT" "Employer"

# --- Subscriptions: append new transaction rows ---
$ws = $wb.Worksheets.Item("Subscriptions")
Set-TxnRow $ws 18 "09/06/2023" 37 "INSTANTLY.AI SHERIDAN US  09/05/23%% Card 15 #9560" "Instantly.ai"
Set-TxnRow $ws 19 "08/06/2023" 37 "INSTANTLY.AI SHERIDAN US  08/05/23%% Card 15 #9560" "Instantly.ai"
Set-TxnRow $ws 20 "07/06/2023" 37 "INSTANTLY.AI SHERIDAN US  07/05/23%% Card 15 #9560" "Instantly.ai"

# --- Uncertain Expenses: append new transaction rows ---
$ws = $wb.Worksheets.Item("Uncertain Expenses")
Set-TxnRow $ws 20 "12/13/2023" 10 "To Share 00 REF# 30479632" "Bank"
Set-TxnRow $ws 21 "11/29/2023" 15 "ODP Fee %% ACH ECC WEB %% ACH Trace 091000016840533" "Bank"
Set-TxnRow $ws 22 "11/29/2023" 210 "From Share 00" "Bank"
Set-TxnRow $ws 23 "12/13/2023" 10 "To Share 00 REF# 30479632" "Bank"

# --- Weekly Budget: insert new header-ish "Week Start" row at row 2, ---
# --- shifting all weekly data down by one, then correct a handful of ---
# --- weekly income/expense totals affected by the new transactions ---
$ws = $wb.Worksheets.Item("Weekly Budget")
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = "Week Start"
$ws.Cells.Item(2, 2).Value = "Income"
$ws.Cells.Item(2, 3).Value = "Expenses"
$ws.Cells.Item(2, 4).Value = "Balance"

$ws.Cells.Item(12, 2).Value = 4058
$ws.Cells.Item(12, 3).Value = 16345
$ws.Cells.Item(12, 4).Value = -12287
$ws.Cells.Item(17, 2).Value = 140
$ws.Cells.Item(17, 3).Value = 5021
$ws.Cells.Item(17, 4).Value = -4881
$ws.Cells.Item(21, 2).Value = 215
$ws.Cells.Item(21, 3).Value = 889
$ws.Cells.Item(21, 4).Value = -674
$ws.Cells.Item(43, 2).Value = 747
$ws.Cells.Item(43, 3).Value = 3752
$ws.Cells.Item(43, 4).Value = -3005
$ws.Cells.Item(44, 2).Value = 4358
$ws.Cells.Item(44, 3).Value = 7766
$ws.Cells.Item(44, 4).Value = -3408
$ws.Cells.Item(46, 2).Value = 600
$ws.Cells.Item(46, 3).Value = 4424
$ws.Cells.Item(46, 4).Value = -3824
$ws.Cells.Item(49, 2).Value = 704
$ws.Cells.Item(49, 3).Value = 6786
$ws.Cells.Item(49, 4).Value = -6082
$ws.Cells.Item(50, 2).Value = 6120
$ws.Cells.Item(50, 3).Value = 18942
$ws.Cells.Item(50, 4).Value = -12822
$ws.Cells.Item(52, 2).Value = 16274
$ws.Cells.Item(52, 3).Value = 9295
$ws.Cells.Item(52, 4).Value = 6979
$ws.Cells.Item(53, 2).Value = 304
$ws.Cells.Item(53, 3).Value = 5482
$ws.Cells.Item(53, 4).Value = -5178
$ws.Cells.Item(54, 2).Value = 4678
$ws.Cells.Item(54, 3).Value = 8714
$ws.Cells.Item(54, 4).Value = -4036

# --- Balances: fix bank account label casing + total amount ---
$ws = $wb.Worksheets.Item("Balances")
$ws.Cells.Item(3, 1).Value = "bank"
$ws.Cells.Item(3, 2).Value = 3000000

# --- Balance Summary: refresh totals and insert a new "Total Balance" row ---
$ws = $wb.Worksheets.Item("Balance Summary")
$ws.Cells.Item(2, 2).Value = 281255.97
$ws.Cells.Item(3, 2).Value = 1399055.88
$ws.Cells.Item(4, 2).Value = -1117799.91
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).Value = "Total Balance"
$ws.Cells.Item(5, 2).Value = 3000000
$ws.Cells.Item(6, 2).Value = 1882200.09
$ws.Cells.Item(7, 2).Value = 5156.712575342466
$ws.Cells.Item(8, 2).Value = 36196.15557692308
$ws.Cells.Item(9, 2).Value = 1882200.09
